$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Participant ID" Cypher query in cell B2 with the new
# query text (uses apoc.coll.sort instead of ORDER BY ... LIMIT 100,
# and resolves the study/participant through optional matches on the
# genomic_info node).
$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH p, g
WHERE g.instrument_model in ['HiSeq X Five']
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
'@

$ws.Range("B2").Value = $newParticipantQuery

# Move the active selection from E4 to C2.
$ws.Range("C2").Select()
